# 55208.xlsx: adapt column header formatting to respective input file names
#  - "<field>_old" -> "<field>_FV2410"
#  - "<field>_new" -> "<field>_FV2504"
#  - wrap the A1:U57 range in an Excel Table ("Table1") so the new header
#    names are also carried by the table's column definitions
#  - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row (row 1, columns A:U) -----------------------
# "<name>_old" -> "<name>_FV2410" and "<name>_new" -> "<name>_FV2504"
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($header -like "*_old") {
        $cell.Value2 = ($header -replace "_old$", "_FV2410")
    } elseif ($header -like "*_new") {
        $cell.Value2 = ($header -replace "_new$", "_FV2504")
    }
}

# --- Turn the header+data range into a native Excel Table --------------
# This creates xl/tables/table1.xml (with the already-renamed column
# headers) and wires up the worksheet's <tableParts>.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ---------------------------------------------
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
